$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "최종점수" (K) column values
$ws.Range("K2").Value = 49.1
$ws.Range("K3").Value = 48.3
$ws.Range("K4").Value = 47.1
$ws.Range("K5").Value = 41.5
$ws.Range("K6").Value = 40.3

# Update "MACRO_SCORE" (N) column values
$ws.Range("N2").Value = 53.62998959737769
$ws.Range("N3").Value = 53.62998959737769
$ws.Range("N4").Value = 53.62998959737769
$ws.Range("N5").Value = 53.62998959737769
$ws.Range("N6").Value = 53.62998959737769
